$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add extra columns to show info on 2 bachelor levels.
# Each of the two "в тому числі" detail groups (яка відноситься до "випуск"
# та до "прийом") gets 2 new trailing sub-columns: "молодших бакалаврів" and
# "молодших фахових бакалаврів".
# ---------------------------------------------------------------------------

# Insert the second pair of columns first (further right), so the column
# letters used for the first insertion point are unaffected.
$ws.Range("N1:O1").EntireColumn.Insert()
$ws.Range("I1:J1").EntireColumn.Insert()

# Re-establish the merged header cells that should now span the 2 extra
# columns in each detail group.
$ws.Range("D5:H5").UnMerge()
$ws.Range("D5:J5").Merge()

$ws.Range("K5:O5").UnMerge()
$ws.Range("K5:Q5").Merge()

$ws.Range("E6:H6").UnMerge()
$ws.Range("E6:J6").Merge()

$ws.Range("L6:O6").UnMerge()
$ws.Range("L6:Q6").Merge()

# New leaf header labels for row 7.
$ws.Range("I7").Value = "молодших бакалаврів"
$ws.Range("J7").Value = "молодших фахових бакалаврів"
$ws.Range("P7").Value = "молодших бакалаврів"
$ws.Range("Q7").Value = "молодших фахових бакалаврів"

# Match formatting of the new cells to their neighbours.
$ws.Range("H7").Copy()
$ws.Range("I7:J7").PasteSpecial(-4122)
$ws.Range("O7").Copy()
$ws.Range("P7:Q7").PasteSpecial(-4122)

$ws.Range("H4:H9").Copy()
$ws.Range("I4:J9").PasteSpecial(-4122)
$ws.Range("O4:O9").Copy()
$ws.Range("P4:Q9").PasteSpecial(-4122)

$excel.CutCopyMode = 0
